# Scheduled runner update: refresh computed leve-profit figures across sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to reflect latest market prices.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1668402.6
$ws.Range("J17").Value = 1668402.6
$ws.Range("L17").Value = 5005207.800000001
$ws.Range("N17").Value = -5005543.800000001

$ws.Range("H39").Value = 432
$ws.Range("I39").Value = 338.4
$ws.Range("K39").Value = 1015.2
$ws.Range("M39").Value = -719.1999999999999

$ws.Range("H137").Value = 23409.46
$ws.Range("I137").Value = 27244.117
$ws.Range("K137").Value = 81732.351
$ws.Range("M137").Value = -79182.351

$ws.Range("H138").Value = 4010.6985
$ws.Range("J138").Value = 4080
$ws.Range("L138").Value = 12240
$ws.Range("N138").Value = -22520

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H34").Value = 50674
$ws.Range("I34").Value = 50674
$ws.Range("K34").Value = 50674
$ws.Range("M34").Value = -50403

$ws.Range("H45").Value = 2124.5
$ws.Range("J45").Value = 1249
$ws.Range("L45").Value = 1249
$ws.Range("N45").Value = -2003

$ws.Range("H97").Value = 1223.9524
$ws.Range("I97").Value = 447.46667
$ws.Range("K97").Value = 447.46667
$ws.Range("M97").Value = 48.53332999999998

$ws.Range("H102").Value = 1780.4375
$ws.Range("I102").Value = 1642.3846
$ws.Range("K102").Value = 1642.3846
$ws.Range("M102").Value = -20.38460000000009

$ws.Range("H122").Value = 2079.4827
$ws.Range("I122").Value = 2085.7778
$ws.Range("K122").Value = 6257.3334
$ws.Range("M122").Value = -3807.3334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2794.5
$ws.Range("I99").Value = 2983.6206
$ws.Range("J99").Value = 966.3333
$ws.Range("K99").Value = 2983.6206
$ws.Range("L99").Value = 966.3333
$ws.Range("M99").Value = -1485.6206
$ws.Range("N99").Value = -3962.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 246803.56
$ws.Range("I31").Value = 304875.28
$ws.Range("J31").Value = 7257.75
$ws.Range("K31").Value = 304875.28
$ws.Range("L31").Value = 7257.75
$ws.Range("M31").Value = -304580.28
$ws.Range("N31").Value = -7847.75

$ws.Range("H34").Value = 246803.56
$ws.Range("I34").Value = 304875.28
$ws.Range("J34").Value = 7257.75
$ws.Range("K34").Value = 304875.28
$ws.Range("L34").Value = 7257.75
$ws.Range("M34").Value = -304673.28
$ws.Range("N34").Value = -7661.75

$ws.Range("H64").Value = 51423.668
$ws.Range("J64").Value = 51423.668
$ws.Range("L64").Value = 51423.668
$ws.Range("N64").Value = -51919.668

$ws.Range("H67").Value = 51423.668
$ws.Range("J67").Value = 51423.668
$ws.Range("L67").Value = 51423.668
$ws.Range("N67").Value = -53139.668

$ws.Range("H86").Value = 3889
$ws.Range("I86").Value = 2692.5557
$ws.Range("K86").Value = 2692.5557
$ws.Range("M86").Value = -1569.5557

$ws.Range("H89").Value = 3889
$ws.Range("I89").Value = 2692.5557
$ws.Range("K89").Value = 13462.7785
$ws.Range("M89").Value = -7846.7785

$ws.Range("H107").Value = 6402.125
$ws.Range("I107").Value = 1542
$ws.Range("K107").Value = 1542
$ws.Range("M107").Value = 378

$ws.Range("H122").Value = 3816.0908
$ws.Range("I122").Value = 2997.4443
$ws.Range("J122").Value = 7500
$ws.Range("K122").Value = 8992.332900000001
$ws.Range("L122").Value = 22500
$ws.Range("M122").Value = -6542.332900000001
$ws.Range("N122").Value = -27400

$ws.Range("H132").Value = 3502.4443
$ws.Range("I132").Value = 2022.3077
$ws.Range("K132").Value = 6066.9231
$ws.Range("M132").Value = -3536.9231

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 5499.5
$ws.Range("I5").Value = 10000
$ws.Range("J5").Value = 999
$ws.Range("K5").Value = 30000
$ws.Range("L5").Value = 2997
$ws.Range("M5").Value = -29888
$ws.Range("N5").Value = -3221

$ws.Range("H68").Value = 13892264
$ws.Range("J68").Value = 5758.3335
$ws.Range("L68").Value = 17275.0005
$ws.Range("N68").Value = -18897.0005

$ws.Range("H71").Value = 13892264
$ws.Range("J71").Value = 5758.3335
$ws.Range("L71").Value = 51825.0015
$ws.Range("N71").Value = -59937.0015

$ws.Range("H107").Value = 384.39285
$ws.Range("I107").Value = 337.1111
$ws.Range("J107").Value = 469.5
$ws.Range("K107").Value = 1011.3333
$ws.Range("L107").Value = 1408.5
$ws.Range("M107").Value = 908.6667
$ws.Range("N107").Value = -5248.5

$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").ClearContents()
$ws.Range("N120").Value = 0

$ws.Range("H132").Value = 6919.0435
$ws.Range("I132").Value = 9867.429
$ws.Range("J132").Value = 2332.6667
$ws.Range("K132").Value = 88806.861
$ws.Range("L132").Value = 20994.0003
$ws.Range("M132").Value = -86276.861
$ws.Range("N132").Value = -26054.0003

$ws.Range("H135").Value = 5499.5
$ws.Range("I135").Value = 10000
$ws.Range("J135").Value = 999
$ws.Range("K135").Value = 90000
$ws.Range("L135").Value = 8991
$ws.Range("M135").Value = -87465
$ws.Range("N135").Value = -14061

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 3182.889
$ws.Range("I97").Value = 1893.4
$ws.Range("K97").Value = 1893.4
$ws.Range("M97").Value = -1397.4

$ws.Range("H102").Value = 3656.0908
$ws.Range("I102").Value = 4008.0557
$ws.Range("K102").Value = 4008.0557
$ws.Range("M102").Value = -2386.0557

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H29").Value = 40000
$ws.Range("I29").Value = 30000
$ws.Range("K29").Value = 30000
$ws.Range("M29").Value = -29705

$ws.Range("H38").Value = 52749.625
$ws.Range("J38").Value = 52749.625
$ws.Range("L38").Value = 52749.625
$ws.Range("N38").Value = -53569.625

$ws.Range("H122").Value = 6790
$ws.Range("I122").Value = 6731.1113
$ws.Range("K122").Value = 20193.3339
$ws.Range("M122").Value = -17743.3339

$ws.Range("H136").Value = 2165.973
$ws.Range("I136").Value = 1999.96
$ws.Range("K136").Value = 5999.88
$ws.Range("M136").Value = -3449.88

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H34").Value = 26675
$ws.Range("I34").Value = 26675
$ws.Range("K34").Value = 26675
$ws.Range("M34").Value = -26472

$ws.Range("H37").Value = 56605.668
$ws.Range("I37").Value = 39820
$ws.Range("J37").Value = 64998.5
$ws.Range("K37").Value = 39820
$ws.Range("L37").Value = 64998.5
$ws.Range("M37").Value = -39617
$ws.Range("N37").Value = -65404.5

$ws.Range("H92").Value = 30000
$ws.Range("J92").Value = 30000
$ws.Range("L92").Value = 30000
$ws.Range("N92").Value = -34992

$ws.Range("H107").Value = 1199.5454
$ws.Range("I107").Value = 1047.2222
$ws.Range("K107").Value = 3141.6666
$ws.Range("M107").Value = -1221.6666

$ws.Range("H113").Value = 1161.7097
$ws.Range("I113").Value = 935.6667
$ws.Range("K113").Value = 2807.0001
$ws.Range("M113").Value = -637.0001000000002

$ws.Range("H132").Value = 2070.6287
$ws.Range("I132").Value = 1560.9375
$ws.Range("K132").Value = 4682.8125
$ws.Range("M132").Value = -2152.8125
